# Auto-generated edit script applying the diff to Asura_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 13
$ws.Range("H13").Value = 34999.5
$ws.Range("J13").Value = 19999
$ws.Range("L13").Value = 19999
$ws.Range("N13").Value = -20337

# ALC!row 116
$ws.Range("H116").Value = 15386838
$ws.Range("I116").Value = 66668300
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 66668300
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = -66664858
$ws.Range("N116").Value = -9284

# ALC!row 129
$ws.Range("H129").Value = 1241.3513
$ws.Range("J129").Value = 1396.129
$ws.Range("L129").Value = 4188.387
$ws.Range("N129").Value = -14188.387

# ALC!row 132
$ws.Range("H132").Value = 1860.2903
$ws.Range("I132").Value = 1802.5555
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 5407.666499999999
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -2877.666499999999
$ws.Range("N132").Value = -11810

# ALC!row 138
$ws.Range("H138").Value = 3997.9622
$ws.Range("I138").Value = 2261.8
$ws.Range("J138").Value = 4683.2896
$ws.Range("K138").Value = 6785.400000000001
$ws.Range("L138").Value = 14049.8688
$ws.Range("M138").Value = -1645.400000000001
$ws.Range("N138").Value = -24329.8688

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 32
$ws.Range("H32").Value = 11507.591
$ws.Range("I32").Value = 12500.741
$ws.Range("J32").Value = 4307.25
$ws.Range("K32").Value = 12500.741
$ws.Range("L32").Value = 4307.25
$ws.Range("M32").Value = -12213.741
$ws.Range("N32").Value = -4881.25

# ARM!row 121
$ws.Range("H121").Value = 56660
$ws.Range("J121").Value = 56660
$ws.Range("L121").Value = 56660
$ws.Range("N121").Value = -60154

# ARM!row 122
$ws.Range("H122").Value = 2181.818
$ws.Range("I122").Value = 2284.7827
$ws.Range("J122").Value = 1945
$ws.Range("K122").Value = 6854.348100000001
$ws.Range("L122").Value = 5835
$ws.Range("M122").Value = -4404.348100000001
$ws.Range("N122").Value = -10735

# ARM!row 125
$ws.Range("H125").Value = 49999
$ws.Range("J125").Value = 49999
$ws.Range("L125").Value = 49999
$ws.Range("N125").Value = -59839

# ARM!row 131
$ws.Range("H131").Value = 70000
$ws.Range("J131").Value = 70000
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 6
$ws.Range("H6").Value = 1847.6316
$ws.Range("J6").Value = 1988.4706
$ws.Range("L6").Value = 5965.4118
$ws.Range("N6").Value = -6191.4118

# CUL!row 69
$ws.Range("H69").Value = 700
$ws.Range("I69").Value = 700
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2100
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("M69").Value = -1289

# CUL!row 72
$ws.Range("H72").Value = 700
$ws.Range("I72").Value = 700
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 6300
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("M72").Value = -2244

# CUL!row 92
$ws.Range("H92").Value = 1312
$ws.Range("I92").Value = 748
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 2244
$ws.Range("L92").Value = 4500
$ws.Range("M92").Value = -996
$ws.Range("N92").Value = -6996

# CUL!row 131
$ws.Range("H131").Value = 13890099
$ws.Range("J131").Value = 14286909
$ws.Range("L131").Value = 42860727
$ws.Range("N131").Value = -42870807

# CUL!row 137
$ws.Range("H137").Value = 2646.6667
$ws.Range("I137").Value = 931
$ws.Range("J137").Value = 4791.25
$ws.Range("K137").Value = 2793
$ws.Range("L137").Value = 14373.75
$ws.Range("M137").Value = 2307
$ws.Range("N137").Value = -24573.75

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 31
$ws.Range("H31").Value = 20632.75
$ws.Range("I31").Value = 20632.75
$ws.Range("K31").Value = 20632.75
$ws.Range("M31").Value = -20340.75

# GSM!row 37
$ws.Range("H37").Value = 20632.75
$ws.Range("I37").Value = 20632.75
$ws.Range("K37").Value = 20632.75
$ws.Range("M37").Value = -20355.75

# GSM!row 109
$ws.Range("H109").Value = 20284.334
$ws.Range("J109").Value = 20284.334
$ws.Range("L109").Value = 20284.334
$ws.Range("N109").Value = -22364.334

# GSM!row 123
$ws.Range("H123").Value = 8925
$ws.Range("J123").Value = 8925
$ws.Range("L123").Value = 8925
$ws.Range("N123").Value = -13825

# GSM!row 131
$ws.Range("H131").Value = 37654
$ws.Range("J131").Value = 37654
$ws.Range("L131").Value = 37654
$ws.Range("N131").Value = -47734

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# LTW!row 29
$ws.Range("H29").Value = 14000
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 23000
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 23000
$ws.Range("M29").Value = -4705
$ws.Range("N29").Value = -23590

# LTW!row 40
$ws.Range("H40").Value = 6111.2856
$ws.Range("I40").Value = 9640
$ws.Range("J40").Value = 4699.8
$ws.Range("K40").Value = 9640
$ws.Range("L40").Value = 4699.8
$ws.Range("M40").Value = -9504
$ws.Range("N40").Value = -4971.8

# LTW!row 46
$ws.Range("H46").Value = 1533.3334
$ws.Range("I46").Value = 1250
$ws.Range("J46").Value = 1675
$ws.Range("K46").Value = 1250
$ws.Range("L46").Value = 1675
$ws.Range("M46").Value = -1062
$ws.Range("N46").Value = -2051

# LTW!row 122
$ws.Range("H122").Value = 16670262
$ws.Range("I122").Value = 3707.3333
$ws.Range("J122").Value = 33336816
$ws.Range("K122").Value = 11121.9999
$ws.Range("L122").Value = 100010448
$ws.Range("M122").Value = -8671.999899999999
$ws.Range("N122").Value = -100015348

# LTW!row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 94
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802

# WVR!row 123
$ws.Range("H123").Value = 32852.758
$ws.Range("J123").Value = 32852.758
$ws.Range("L123").Value = 32852.758
$ws.Range("N123").Value = -42652.758

# WVR!row 125
$ws.Range("H125").Value = 60643.332
$ws.Range("J125").Value = 60643.332
$ws.Range("L125").Value = 60643.332
$ws.Range("N125").Value = -70483.33199999999

# WVR!row 126
$ws.Range("H126").Value = 11689.1875
$ws.Range("I126").Value = 18753.5
$ws.Range("J126").Value = 4624.875
$ws.Range("K126").Value = 56260.5
$ws.Range("L126").Value = 13874.625
$ws.Range("M126").Value = -53790.5
$ws.Range("N126").Value = -18814.625

# WVR!row 132
$ws.Range("H132").Value = 2074.28
$ws.Range("I132").Value = 1327.9656
$ws.Range("J132").Value = 3104.9048
$ws.Range("K132").Value = 3983.8968
$ws.Range("L132").Value = 9314.714399999999
$ws.Range("M132").Value = -1453.8968
$ws.Range("N132").Value = -14374.7144
